$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.594.71"
$ws.Range("E2").Value = "  +1.67%  "

$ws.Range("D3").Value = "3.162.00"
$ws.Range("E3").Value = "  +1.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.31"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.47"
$ws.Range("E6").Value = "  +1.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +14.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.29"
$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("E10").Value = "  +6.10%  "

$ws.Range("E11").Value = "  +4.11%  "

$ws.Range("E12").Value = "  +2.50%  "

$ws.Range("D13").Value = "3.708.29"
$ws.Range("E13").Value = "  +1.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.73"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000170"
$ws.Range("E15").Value = "  +3.49%  "

$ws.Range("D16").Value = "58.646.08"
$ws.Range("E16").Value = "  +1.53%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.188.60"
$ws.Range("E17").Value = "  +2.35%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.24"
$ws.Range("E18").Value = "  +3.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  +2.16%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "376.10"
$ws.Range("E20").Value = "  +3.96%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.10"
$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.529"
$ws.Range("E23").Value = "  +4.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.73"
$ws.Range("E24").Value = "  +1.08%  "

$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.27"
$ws.Range("E27").Value = "  +13.11%  "

$ws.Range("D28").Value = "0.0₃0863"
$ws.Range("E28").Value = "  -0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.38"
$ws.Range("E29").Value = "  +4.66%  "

$ws.Range("E30").Value = "  +0.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.03"
$ws.Range("E31").Value = "  -1.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.16"
$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("E33").Value = "  +0.85%  "

$ws.Range("E34").Value = "  +4.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.62"
$ws.Range("E35").Value = "  -1.90%  "

$ws.Range("E36").Value = "  +4.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.02"
$ws.Range("E37").Value = "  -2.03%  "

$ws.Range("D38").Value = "2.683.46"
$ws.Range("E38").Value = "  +7.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0696"
$ws.Range("E39").Value = "  +3.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.68"
$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("E41").Value = "  +6.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.721"
$ws.Range("E42").Value = "  +3.62%  "

$ws.Range("E43").Value = "  +3.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0289"
$ws.Range("E44").Value = "  +7.04%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "3.204.74"
$ws.Range("E46").Value = "  +1.51%  "

$ws.Range("E47").Value = "  +14.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.22"
$ws.Range("E48").Value = "  +2.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.977"
$ws.Range("E49").Value = "  -0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.03"
$ws.Range("E50").Value = "  +1.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.747"
$ws.Range("E51").Value = "  +1.54%  "

# Reset style to Normal so no extra numFmt/style remains applied to the cell
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
